# "Arrumei a forma como pegar o codigo e o preco medio"
# Fill in the Acura model/year/FIPE-code/average-price table: extend the
# existing Integra GS 1.8 rows with their FIPE code + prices, add the
# Legend 3.2/3.5 and NSX 3.0 model rows, each with one row per model year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Route the numeric-looking string through a formula-text scratch cell and
    # paste-special (values only) onto the target so Excel stores it as a plain
    # shared-string cell instead of auto-converting it to a number.
    $scratch = $ws.Range("Z1")
    $scratch.Formula = "=""" + $val + """"
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

$ws.Range("A2").Value = "Acura"
$ws.Range("B2").Value = "Integra GS 1.8"
$ws.Range("C2").Value = "1992 Gasolina"
$ws.Range("D2").Value = "038003-2"
Set-TextValue "E2" " 11097.00"

$ws.Range("A3").Value = "Acura"
$ws.Range("B3").Value = "Integra GS 1.8"
$ws.Range("C3").Value = "1991 Gasolina"
$ws.Range("D3").Value = "038003-2"
Set-TextValue "E3" " 10366.00"

$ws.Range("A4").Value = "Acura"
$ws.Range("B4").Value = "Legend 3.2/3.5"
$ws.Range("C4").Value = "1998 Gasolina"
$ws.Range("D4").Value = "038002-4"
Set-TextValue "E4" " 25397.00"

$ws.Range("A5").Value = "Acura"
$ws.Range("B5").Value = "Legend 3.2/3.5"
$ws.Range("C5").Value = "1997 Gasolina"
$ws.Range("D5").Value = "038002-4"
Set-TextValue "E5" " 22580.00"

$ws.Range("A6").Value = "Acura"
$ws.Range("B6").Value = "Legend 3.2/3.5"
$ws.Range("C6").Value = "1996 Gasolina"
$ws.Range("D6").Value = "038002-4"
Set-TextValue "E6" " 19084.00"

$ws.Range("A7").Value = "Acura"
$ws.Range("B7").Value = "Legend 3.2/3.5"
$ws.Range("C7").Value = "1995 Gasolina"
$ws.Range("D7").Value = "038002-4"
Set-TextValue "E7" " 14802.00"

$ws.Range("A8").Value = "Acura"
$ws.Range("B8").Value = "Legend 3.2/3.5"
$ws.Range("C8").Value = "1994 Gasolina"
$ws.Range("D8").Value = "038002-4"
Set-TextValue "E8" " 14219.00"

$ws.Range("A9").Value = "Acura"
$ws.Range("B9").Value = "Legend 3.2/3.5"
$ws.Range("C9").Value = "1993 Gasolina"
$ws.Range("D9").Value = "038002-4"
Set-TextValue "E9" " 14219.00"

$ws.Range("A10").Value = "Acura"
$ws.Range("B10").Value = "Legend 3.2/3.5"
$ws.Range("C10").Value = "1992 Gasolina"
$ws.Range("D10").Value = "038002-4"
Set-TextValue "E10" " 14219.00"

$ws.Range("A11").Value = "Acura"
$ws.Range("B11").Value = "Legend 3.2/3.5"
$ws.Range("C11").Value = "1991 Gasolina"
$ws.Range("D11").Value = "038002-4"
Set-TextValue "E11" " 14219.00"

$ws.Range("A12").Value = "Acura"
$ws.Range("B12").Value = "NSX 3.0"
$ws.Range("C12").Value = "1995 Gasolina"
$ws.Range("D12").Value = "038001-6"
Set-TextValue "E12" " 40991.00"

$ws.Range("A13").Value = "Acura"
$ws.Range("B13").Value = "NSX 3.0"
$ws.Range("C13").Value = "1994 Gasolina"
$ws.Range("D13").Value = "038001-6"
Set-TextValue "E13" " 39550.00"

$ws.Range("A14").Value = "Acura"
$ws.Range("B14").Value = "NSX 3.0"
$ws.Range("C14").Value = "1993 Gasolina"
$ws.Range("D14").Value = "038001-6"
Set-TextValue "E14" " 36538.00"

$ws.Range("A15").Value = "Acura"
$ws.Range("B15").Value = "NSX 3.0"
$ws.Range("C15").Value = "1992 Gasolina"
$ws.Range("D15").Value = "038001-6"
Set-TextValue "E15" " 33397.00"

$ws.Range("A16").Value = "Acura"
$ws.Range("B16").Value = "NSX 3.0"
$ws.Range("C16").Value = "1991 Gasolina"
$ws.Range("D16").Value = "038001-6"
Set-TextValue "E16" " 33397.00"
